$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "29.902.26"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +0.04%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.901.45"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +0.28%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.000"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  -0.02%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "0.7994"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +5.86%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "241.02"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +0.39%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "1.000"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -0.07%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3117"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +2.40%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "26.16"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +2.80%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.06877"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +0.61%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07971"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +0.00%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.899.04"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -0.02%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.7349"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -1.52%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "5.167"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -0.53%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "92.34"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +1.32%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "29.897.34"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -0.01%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "13.91"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -0.12%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "5.853"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -1.57%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "244.32"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +0.62%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.000007702"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -0.35%  "
$ws.Range("E21").Value = "  -0.06%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "2.151.87"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -0.32%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "1.001"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +0.02%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "6.914"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -0.59%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "167.41"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +1.24%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "9.170"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -0.64%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.1423"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +9.56%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "18.81"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +0.39%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.031"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +0.13%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.355"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -3.30%  "
$ws.Range("E31").Value = "  -0.08%  "
$ws.Range("E32").Value = "  +0.19%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.05592"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +4.42%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "4.055"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +0.83%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.257"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +0.60%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.7278"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +0.32%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.717"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -0.03%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.01929"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +0.79%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "2.783"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -0.19%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.4396"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -0.20%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "5.992"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -2.93%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "71.89"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -0.54%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.9998"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -0.07%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.8359"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +1.44%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "1.859"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -2.35%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "100.43"
$ws.Range("D46").Style = "Normal"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "7.549"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -0.15%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "9.706"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -0.24%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "976.17"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +7.68%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "2.058.22"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +0.00%  "
$ws.Range("E51").Value = "  -0.17%  "
